# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 00:49"

# --- Swap country labels that changed ranking order ---
# Colombia overtakes España (rows 8/9)
$ws.Range("A8").Value = "Colombia"
$ws.Range("A9").Value = "España"

# Bulgaria overtakes Tunez (rows 83/84)
$ws.Range("A83").Value = "Bulgaria"
$ws.Range("A84").Value = "Tunez"

# Nueva Caledonia overtakes Santa Lucia (rows 207/208) - values identical, label swap only
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"

# --- Update numeric data (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) for the affected rows ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7715070
$ws.Range("C4").Value = 35984
$ws.Range("D4").Value = 4925310
$ws.Range("E4").Value = 2574068
$ws.Range("G4").Value = 660
$ws.Range("H4").Value = 215692

# Row 6: Brasil
$ws.Range("D6").Value = 4352871
$ws.Range("E6").Value = 468776

# Row 8: now Colombia (new totals, overtakes España)
$ws.Range("B8").Value = 869808
$ws.Range("C8").Value = 7650
$ws.Range("D8").Value = 770812
$ws.Range("E8").Value = 71979
$ws.Range("G8").Value = 173
$ws.Range("H8").Value = 27017

# Row 9: now España (previous row-8 values, unchanged since last pull)
$ws.Range("B9").Value = 865631
$ws.Range("C9").Value = 12793
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 261
$ws.Range("H9").Value = 32486

# Row 29: Canada
$ws.Range("B29").Value = 171323
$ws.Range("C29").Value = 2363
$ws.Range("D29").Value = 143994
$ws.Range("E29").Value = 17799
$ws.Range("G29").Value = 26
$ws.Range("H29").Value = 9530

# Row 42: Egipto
$ws.Range("B42").Value = 103902
$ws.Range("C42").Value = 121
$ws.Range("D42").Value = 97449
$ws.Range("E42").Value = 452
$ws.Range("G42").Value = 11
$ws.Range("H42").Value = 6001

# Row 58: Nigeria
$ws.Range("B58").Value = 59583
$ws.Range("C58").Value = 118
$ws.Range("D58").Value = 51308
$ws.Range("E58").Value = 7162

# Row 83: now Bulgaria (new totals, overtakes Tunez)
$ws.Range("B83").Value = 22306
$ws.Range("C83").Value = 436
$ws.Range("D83").Value = 15310
$ws.Range("E83").Value = 6134
$ws.Range("G83").Value = 8
$ws.Range("H83").Value = 862

# Row 84: now Tunez (previous row-83 values, unchanged since last pull)
$ws.Range("B84").Value = 22230
$ws.Range("D84").Value = 5032
$ws.Range("E84").Value = 16877
$ws.Range("H84").Value = 321

# Row 115: Mauritania
$ws.Range("B115").Value = 7529
$ws.Range("C115").Value = 6
$ws.Range("D115").Value = 7208
$ws.Range("E115").Value = 159

# Row 133: Trinidad yTobago
$ws.Range("B133").Value = 4846
$ws.Range("C133").Value = 79
$ws.Range("E133").Value = 1812

# Row 198: Islas Virgenes Britanicas
$ws.Range("D198").Value = 69
$ws.Range("E198").Value = 1
